$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 stays the text label "HK_R_acc_G" but now refers to a newly added
# duplicate shared-string entry (sharedStrings count/uniqueCount 38 -> 39).
$ws.Range("A1").Value = "HK_R_acc_G"

$values = @(
    70.782208588957047,
    68.788343558282207,
    65.874233128834362,
    65.030674846625772,
    69.785276073619627,
    72.315950920245399,
    69.401840490797554,
    69.401840490797554,
    69.401840490797554,
    69.478527607361968,
    71.395705521472394,
    71.165644171779135,
    70.475460122699388,
    70.552147239263803,
    70.475460122699388,
    69.248466257668724,
    71.779141104294482,
    71.779141104294482,
    70.168711656441715,
    61.579754601227002,
    70.245398773006144,
    65.874233128834362,
    66.257668711656436,
    65.950920245398777,
    69.708588957055213,
    69.708588957055213,
    69.401840490797554,
    69.708588957055213,
    69.785276073619627,
    70.245398773006144,
    62.116564417177912,
    66.641104294478524,
    66.871165644171782,
    63.880368098159515,
    63.957055214723923,
    71.395705521472394,
    63.880368098159515,
    63.803680981595093,
    65.950920245398777,
    68.941717791411037,
    69.171779141104295,
    68.865030674846622,
    68.251533742331276,
    68.481595092024534,
    67.791411042944787,
    68.481595092024534,
    64.647239263803684,
    68.634969325153378,
    68.558282208588963
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
